# Colour the "2 quintes ou octaves consécutives" bullet the same green
# (RGB 92D050) already used by the other bullets in this list.

$d = $word.ActiveDocument

$searchText = "Il est interdit de faire 2  quintes ou octaves consécutives même par mouvement contraire"
$greenColor = 5296274   # RGB(0x92, 0xD0, 0x50) -> wdColor value for hex 92D050

$rng = $d.Content
$found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand the found range to the whole paragraph (including the
    # paragraph mark) so the colour is applied the same way Word does
    # when you select an entire bulleted line and set the font colour.
    $para = $rng.Paragraphs(1)
    $paraRange = $para.Range
    $paraRange.Font.Color = $greenColor
    Write-Output "Applied green font colour to the target paragraph."
} else {
    Write-Output "Target paragraph text not found."
}
